$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 5, shifting rows 5-7 down to 6-8.
$ws.Rows.Item(5).Insert()

# Copy formatting (e.g. the date style) from the row that is now row 6 (old row 5)
# into the newly inserted row 5, so the new row keeps consistent formatting.
# Restrict the copy to the used columns (A:R) to avoid touching the whole row.
$ws.Range("A6:R6").Copy()
$ws.Range("A5:R5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 5 with the new record's data.
$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(5, 3).Value = "Coquimbo"
$ws.Cells.Item(5, 4).Value = 44977
$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 100114002
$ws.Cells.Item(5, 7).Value = "Camote"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 16500
$ws.Cells.Item(5, 12).Value = 17000
$ws.Cells.Item(5, 13).Value = 16750
$ws.Cells.Item(5, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(5, 15).Value = "Perú"
$ws.Cells.Item(5, 16).Value = 931
$ws.Cells.Item(5, 17).Value = 18
$ws.Cells.Item(5, 18).Value = "Hortaliza"
